# BV20 Post-Preprocessing & Quality Checks template
# Adds a new "Functional resolution in mm." (EXP.RES) row to the
# Preprocessing block, just above the Temporal high-pass-filter rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 is currently a blank spacer row that separates the "Exclusions
# section above from the "Preprocessing" section (old row 19) below.
# Insert a fresh blank row at (old) row 19 so that the existing spacer
# stays at row 18 (to be filled in below) and a brand-new spacer row
# lands at row 19, pushing the "Preprocessing" block (and everything
# under it) down by one row.
$ws.Rows.Item(19).EntireRow.Insert()

# Fill in the previously-blank row 18 with the new EXP.RES field.
# Field_ID (column E) is set first so the new shared strings are added
# to the table in the same order the original authoring tool used.
$ws.Range("E18").Value = "EXP.RES"
$ws.Range("B18").Value = "Functional resolution in mm."
$ws.Range("D18").Value = "Functional resolution in mm."

# Grow Table1 by one row so the table keeps covering the data block.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E32"))

# Update the window scroll position / selection to match where the
# editor was last working.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C18").Select()
